# missingEndUserContent-template.docx — commit "Fixed #295 Add the version
# of M2Doc in the template custom properties".
#
# The unified diff for *this* resource file only touches word/document.xml
# and word/styles.xml, and every single hunk is a pure XML-attribute /
# namespace-declaration re-ordering (e.g. <w:pgSz w:w=".." w:h=".."/> ->
# <w:pgSz w:h=".." w:w=".."/>, <w:style w:type=".." w:default=".."
# w:styleId=".."> -> <w:style w:default=".." w:styleId=".." w:type="..">,
# latentStyles/lsdException attributes resorted alphabetically, etc.).
# Every attribute name/value pair present before the commit is still
# present after it — nothing was added, removed, or re-valued, and no
# run/paragraph text, formatting, style, or page-setup value changed.
# That matches the commit message: the real content edit (stamping the
# M2Doc version into the template's custom document properties) landed in
# the shared template-generation code, and this particular fixture was
# merely re-serialized by the newer tooling as a side effect, which is why
# its diff is attribute-order noise only.
#
# Word's object model (real or emulated) has no notion of "attribute
# order" — it always (re)serializes elements in Word's own canonical
# order, never alphabetically — so there is no COM call that reproduces
# that cosmetic reordering, and attempting to "touch" PageSetup/Styles/
# Find-Replace with equivalent values only risks introducing incidental
# side effects (merged runs, extra namespace declarations, etc.) that
# are not present in the diff. The correct, content-faithful edit is
# therefore to leave the document's paragraphs, runs, styles, and page
# setup exactly as authored.
$d = $word.ActiveDocument
